$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.633.56'
$ws.Range('E2').Value = '  -1.79%  '

$ws.Range('D3').Value = '3.080.91'
$ws.Range('E3').Value = '  -0.32%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.15'
$ws.Range('E5').Value = '  +0.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.48'
$ws.Range('E6').Value = '  -1.81%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '3.081.33'
$ws.Range('E8').Value = '  -0.27%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.441'
$ws.Range('E9').Value = '  +0.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.15'
$ws.Range('E10').Value = '  -3.16%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.108'
$ws.Range('E11').Value = '  -0.97%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.390'
$ws.Range('E12').Value = '  +2.29%  '

$ws.Range('D13').Value = '3.611.37'
$ws.Range('E13').Value = '  -0.43%  '

$ws.Range('E14').Value = '  +1.91%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.40'
$ws.Range('E15').Value = '  -6.39%  '

$ws.Range('E16').Value = '  -1.53%  '

$ws.Range('D17').Value = '57.671.75'
$ws.Range('E17').Value = '  -1.70%  '

$ws.Range('D18').Value = '3.072.45'
$ws.Range('E18').Value = '  -0.74%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.08'
$ws.Range('E19').Value = '  -2.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.64'
$ws.Range('E20').Value = '  -2.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.94'
$ws.Range('E21').Value = '  -3.31%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '340.29'
$ws.Range('E22').Value = '  +0.67%  '

$ws.Range('E23').Value = '  +0.18%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.510'
$ws.Range('E24').Value = '  +0.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.01'
$ws.Range('E25').Value = '  +1.75%  '

$ws.Range('E26').Value = '  -1.91%  '

$ws.Range('E27').Value = '  -0.09%  '

$ws.Range('D28').Value = '0.0₃0910'
$ws.Range('E28').Value = '  -0.64%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.06%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.36'
$ws.Range('E30').Value = '  -3.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').Value = '  -0.87%  '

$ws.Range('E32').Value = '  +2.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.88'
$ws.Range('E33').Value = '  -0.54%  '

$ws.Range('E34').Value = '  -3.48%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.92'
$ws.Range('E35').Value = '  +2.35%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  -0.39%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.13'
$ws.Range('E37').Value = '  +0.50%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.96'
$ws.Range('E38').Value = '  -4.87%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.25'
$ws.Range('E39').Value = '  -3.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0663'
$ws.Range('E40').Value = '  -2.86%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.58'
$ws.Range('E41').Value = '  +11.18%  '

$ws.Range('E42').Value = '  +0.73%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.682'
$ws.Range('E43').Value = '  +2.71%  '

$ws.Range('D44').Value = '3.122.20'
$ws.Range('E44').Value = '  -0.47%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.84'
$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.05%  '

$ws.Range('D47').Value = '2.278.92'
$ws.Range('E47').Value = '  -0.19%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0262'
$ws.Range('E48').Value = '  +2.15%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.990'
$ws.Range('E49').Value = '  +3.38%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.07'
$ws.Range('E50').Value = '  +1.20%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.49'
$ws.Range('E51').Value = '  -1.83%  '
